# "Generate Report for Archive"
# Update the localization-status report:
#  - Status text changes from "Ready for handoff" to "In Translation"
#    (Overview!E2:F2, and the "Status" column on each per-locale sheet)
#  - The "Status" column is narrower now, so shrink its width on every
#    sheet that has one (Overview columns E & F, and column C on the
#    zh-cn / de-de detail sheets)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 13.4101845877511

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $text = [string]$cell.Text
            if ($text -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Overview sheet: Status columns are E (5) and F (6)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn / de-de detail sheets: Status column is C (3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
